$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 502, shifting existing rows 502:580 down to 503:581
$ws.Rows.Item(502).Insert()

# Populate the newly inserted row 502 with the new weekly record
$ws.Range("A502").Value = 9
$ws.Range("B502").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C502").Value = "Metropolitana"
$ws.Range("D502").Value = 44951
$ws.Range("E502").Value = 13
$ws.Range("F502").Value = "Fruta"
$ws.Range("G502").Value = 100108
$ws.Range("H502").Value = "Tropicales y subtropicales"
$ws.Range("I502").Value = 100108002
$ws.Range("J502").Value = "Mango"
$ws.Range("K502").Value = "Sin especificar"
$ws.Range("L502").Value = "Primera"
$ws.Range("M502").Value = 610
$ws.Range("N502").Value = 6000
$ws.Range("O502").Value = 6500
$ws.Range("P502").Value = 6230
$ws.Range("Q502").Value = "$/bandeja 4 kilos"
$ws.Range("R502").Value = "Perú"
$ws.Range("S502").Value = 1558
$ws.Range("T502").Value = 4
